# Apply targeted cell value updates to Sheet1 per the commit diff.
# Data-driven approach: each entry specifies Row (R), Column (C, 1-based), and new Value (V).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{R=2; C=6; V=2.96},
    @{R=2; C=8; V=2.58},
    @{R=2; C=9; V=2.68},
    @{R=2; C=14; V=4.6},
    @{R=2; C=16; V=2.18},
    @{R=2; C=17; V=1.79},
    @{R=2; C=18; V=1.49},
    @{R=2; C=19; V=2.94},
    @{R=2; C=20; V=1.62},
    @{R=2; C=21; V=2.4},
    @{R=2; C=22; V=1.59},
    @{R=2; C=24; V=18.5},
    @{R=2; C=35; V=1000},
    @{R=2; C=40; V=23},
    @{R=2; C=41; V=19.5},
    @{R=3; C=6; V=4.3},
    @{R=3; C=7; V=6.4},
    @{R=3; C=8; V=1.61},
    @{R=3; C=9; V=1.8},
    @{R=3; C=10; V=4.2},
    @{R=3; C=11; V=5.3},
    @{R=3; C=15; V=1.16},
    @{R=3; C=16; V=2.68},
    @{R=3; C=17; V=1.45},
    @{R=3; C=18; V=1.66},
    @{R=3; C=19; V=2.12},
    @{R=3; C=21; V=2.42},
    @{R=3; C=22; V=2.24},
    @{R=3; C=23; V=1.19},
    @{R=3; C=25; V=15},
    @{R=3; C=26; V=15},
    @{R=3; C=27; V=21},
    @{R=3; C=28; V=30},
    @{R=3; C=29; V=13.5},
    @{R=3; C=31; V=19},
    @{R=3; C=32; V=55},
    @{R=3; C=33; V=25},
    @{R=3; C=34; V=21},
    @{R=3; C=35; V=30},
    @{R=3; C=37; V=65},
    @{R=3; C=38; V=60},
    @{R=3; C=40; V=48},
    @{R=3; C=41; V=7.2},
    @{R=4; C=17; V=1.85},
    @{R=5; C=6; V=1.69},
    @{R=5; C=7; V=1.83},
    @{R=5; C=10; V=3.4},
    @{R=5; C=16; V=1.58},
    @{R=6; C=7; V=10.5},
    @{R=6; C=10; V=4.2},
    @{R=6; C=14; V=3.45},
    @{R=6; C=16; V=1.86},
    @{R=6; C=17; V=1.96},
    @{R=6; C=19; V=3.5},
    @{R=6; C=20; V=2.14},
    @{R=6; C=22; V=2.96},
    @{R=6; C=35; V=48},
    @{R=7; C=9; V=2.94},
    @{R=7; C=10; V=3},
    @{R=7; C=27; V=50},
    @{R=7; C=29; V=7.2},
    @{R=7; C=32; V=970},
    @{R=7; C=34; V=20},
    @{R=7; C=37; V=48},
    @{R=8; C=6; V=5.9},
    @{R=8; C=9; V=1.67},
    @{R=8; C=19; V=2.72},
    @{R=8; C=22; V=2.48},
    @{R=8; C=35; V=40},
    @{R=9; C=6; V=1.2},
    @{R=9; C=7; V=1.3},
    @{R=9; C=8; V=13.5},
    @{R=9; C=9; V=19.5},
    @{R=9; C=10; V=6.4},
    @{R=9; C=11; V=9.199999999999999},
    @{R=9; C=12; V=1.17},
    @{R=9; C=14; V=3.55},
    @{R=9; C=15; V=1.09},
    @{R=9; C=16; V=3.55},
    @{R=9; C=17; V=1.3},
    @{R=9; C=18; V=1.9},
    @{R=9; C=19; V=1.75},
    @{R=9; C=22; V=1.05},
    @{R=9; C=23; V=4.6},
    @{R=10; C=6; V=1.89},
    @{R=10; C=7; V=1.93},
    @{R=10; C=8; V=5.2},
    @{R=10; C=9; V=5.7},
    @{R=10; C=10; V=3.3},
    @{R=10; C=11; V=3.5},
    @{R=10; C=14; V=2.76},
    @{R=10; C=16; V=1.58},
    @{R=10; C=18; V=1.21},
    @{R=10; C=20; V=2.22},
    @{R=10; C=21; V=1.72},
    @{R=10; C=22; V=1.21},
    @{R=10; C=23; V=2.06},
    @{R=10; C=25; V=14},
    @{R=10; C=26; V=44},
    @{R=10; C=27; V=200},
    @{R=10; C=28; V=6.4},
    @{R=10; C=29; V=8},
    @{R=10; C=30; V=24},
    @{R=10; C=32; V=11.5},
    @{R=10; C=33; V=12.5},
    @{R=10; C=34; V=28},
    @{R=10; C=36; V=22},
    @{R=10; C=40; V=22},
    @{R=10; C=41; V=210},
    @{R=11; C=6; V=1.17},
    @{R=11; C=7; V=1.2},
    @{R=11; C=8; V=14},
    @{R=11; C=9; V=19.5},
    @{R=11; C=10; V=9.199999999999999},
    @{R=11; C=11; V=11},
    @{R=11; C=12; V=1.15},
    @{R=11; C=14; V=10},
    @{R=11; C=15; V=1.08},
    @{R=11; C=16; V=4.1},
    @{R=11; C=18; V=2.24},
    @{R=11; C=19; V=1.69},
    @{R=11; C=23; V=6},
    @{R=11; C=24; V=75},
    @{R=11; C=25; V=100},
    @{R=11; C=40; V=2.56},
    @{R=12; C=12; V=1.2},
    @{R=12; C=26; V=23},
    @{R=12; C=27; V=36},
    @{R=13; C=6; V=1.52},
    @{R=13; C=7; V=1.53},
    @{R=13; C=8; V=6.8},
    @{R=13; C=9; V=7},
    @{R=13; C=10; V=5},
    @{R=13; C=11; V=5.1},
    @{R=13; C=18; V=1.66},
    @{R=13; C=20; V=1.73},
    @{R=13; C=21; V=2.3},
    @{R=13; C=23; V=2.9},
    @{R=13; C=27; V=190},
    @{R=13; C=28; V=12},
    @{R=13; C=30; V=25},
    @{R=13; C=38; V=25},
    @{R=13; C=40; V=5.7},
    @{R=14; C=7; V=3.9},
    @{R=14; C=14; V=5.9},
    @{R=14; C=16; V=2.68},
    @{R=14; C=18; V=1.69},
    @{R=14; C=19; V=2.16},
    @{R=14; C=20; V=1.48},
    @{R=14; C=21; V=2.62},
    @{R=14; C=22; V=1.87},
    @{R=14; C=23; V=1.35},
    @{R=14; C=38; V=36},
    @{R=15; C=8; V=3.75},
    @{R=15; C=9; V=4.1},
    @{R=15; C=10; V=3.3},
    @{R=15; C=11; V=3.4},
    @{R=15; C=16; V=1.78},
    @{R=15; C=17; V=2.14},
    @{R=15; C=18; V=1.29},
    @{R=15; C=20; V=1.87},
    @{R=15; C=21; V=2.04},
    @{R=15; C=22; V=1.33},
    @{R=15; C=26; V=34},
    @{R=15; C=29; V=7.8},
    @{R=15; C=31; V=70},
    @{R=15; C=36; V=32},
    @{R=16; C=9; V=4.9},
    @{R=16; C=10; V=3.55},
    @{R=16; C=11; V=3.95},
    @{R=16; C=16; V=1.85},
    @{R=16; C=17; V=1.91},
    @{R=16; C=22; V=1.27},
    @{R=17; C=6; V=1.19},
    @{R=17; C=7; V=1.25},
    @{R=17; C=9; V=19},
    @{R=17; C=10; V=6.6},
    @{R=17; C=11; V=12},
    @{R=17; C=15; V=1.09},
    @{R=17; C=18; V=2.16},
    @{R=17; C=19; V=1.66},
    @{R=17; C=20; V=1.75},
    @{R=18; C=8; V=4.5},
    @{R=18; C=9; V=5.3},
    @{R=18; C=39; V=1000}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.R, $u.C).Value = $u.V
}
